# Update cryptos list data (price/volume columns, and two row reorderings
# where coin name/link/price/volume moved between rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.776.04'
$ws.Range("E2").Value = '  +0.65%  '
# Row 3
$ws.Range("D3").Value = '2.617.13'
$ws.Range("E3").Value = '  -0.55%  '
# Row 4
$ws.Range("E4").Value = '  -0.06%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.14%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.38%  '
# Row 7
$ws.Range("E7").Value = '  -0.01%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.543'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.95%  '
# Row 9
$ws.Range("D9").Value = '2.616.77'
$ws.Range("E9").Value = '  -0.52%  '
# Row 10
$ws.Range("E10").Value = '  +7.26%  '
# Row 11
$ws.Range("E11").Value = '  -0.54%  '
# Row 12
$ws.Range("E12").Value = '  -0.23%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.346'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.79%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.47'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.97%  '
# Row 15
$ws.Range("E15").Value = '  +2.80%  '
# Row 16
$ws.Range("D16").Value = '3.093.35'
$ws.Range("E16").Value = '  -0.59%  '
# Row 17
$ws.Range("D17").Value = '67.664.23'
$ws.Range("E17").Value = '  +0.61%  '
# Row 18
$ws.Range("D18").Value = '2.609.10'
$ws.Range("E18").Value = '  -0.84%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '371.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.17%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.12%  '
# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.58%  '
# Row 22
$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.92%  '
# Row 23
$ws.Range("E23").Value = '  -3.36%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.49%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.77%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.05%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.74%  '
# Row 28
$ws.Range("D28").Value = '2.747.16'
$ws.Range("E28").Value = '  -0.61%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000104'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.43%  '
# Row 30
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.01'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.64%  '
# Row 31
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '581.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.96%  '
# Row 32
$ws.Range("E32").Value = '  -1.36%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.78'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.18%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.70%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '
# Row 36
$ws.Range("E36").Value = '  -2.08%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.50'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.87%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.58%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.13'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.49%  '
# Row 40
$ws.Range("E40").Value = '  +3.25%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.367'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.89%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.28'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.07%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.12%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.46%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.01%  '
# Row 46
$ws.Range("E46").Value = '  -2.02%  '
# Row 47
$ws.Range("D47").Value = '0.0₆0304'
$ws.Range("E47").Value = '  +4.42%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '154.84'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.75%  '
# Row 49
$ws.Range("E49").Value = '  -1.58%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.69'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.53%  '
# Row 51
$ws.Range("E51").Value = '  -1.99%  '
